# "split A B C D Channels" - add a new "num of periods sampled" column (E)
# and a new data row (row 3) on the Feuille3 worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuille3")

# New header cell for column E
$ws.Range("E1").Value = "num of periods sampled"

# Widen the new column E to fit the longer header text
$ws.Columns.Item(5).ColumnWidth = 20.83

# Existing row 2 grows slightly taller once the new row is added
$ws.Rows.Item(2).RowHeight = 13.8

# New row 3: another channel sample (15 MHz), with its derived period and
# number of periods sampled relative to row 2's sample count
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Range("A3").Value = 15
$ws.Range("B3").Formula = "=1000/A3"
$ws.Range("E3").Formula = "=D2/B3"

# Put the selection where the author left it after editing
[void]$ws.Range("D17").Select()
